# Insert two new rows before the current row 25 ("cube_field"/"Field"),
# pushing the existing cube_* rows down by two (25-29 -> 27-31), and
# populate the new rows with the victory/score strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("25:26").Insert()

$ws.Range("A25").Value = "victory_title"
$ws.Range("B25").Value = "MISSION COMPLETE"
$ws.Range("A26").Value = "score"
$ws.Range("B26").Value = "Score:"

# Match the updated view state recorded in the workbook: scrolled so row 10
# is at the top, with A26 selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("A26").Select() | Out-Null
